$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Collect "skip" ranges: spans of text whose enclosing run must NOT be
# touched by the generic "{{" / "}}" delimiter loop below, because they are
# handled specially (either the whole-tag runs "{{desconto}}" /
# "{{desconto_extenso}}", or the one run that reads " }}" including its
# leading space).
# ---------------------------------------------------------------------------
$skipRanges = New-Object System.Collections.ArrayList

function InsideSkipRange($start) {
    foreach ($pair in $skipRanges) {
        if ($start -ge $pair[0] -and $start -lt $pair[1]) {
            return $true
        }
    }
    return $false
}

# ---------------------------------------------------------------------------
# Step 1: the two runs whose full text is the whole tag "{{desconto}}" or
# "{{desconto_extenso}}" already carry <w:b/>; the diff only adds
# <w:color w:val="000000"/> to them (the run is not split).
# ---------------------------------------------------------------------------
foreach ($tag in @("{{desconto}}", "{{desconto_extenso}}")) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $found = $rng.Find.Execute($tag, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = $true
        $rng.Font.Color = 0
        [void]$skipRanges.Add(@($rng.Start, $rng.End))
    }
}

# ---------------------------------------------------------------------------
# Step 2: the single run whose full text is " }}" (leading space kept in the
# same run as the closing delimiter) gets bold + black color applied to the
# whole run (space included), matching the diff exactly.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute(" }}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Font.Bold = $true
    $rng.Font.Color = 0
    [void]$skipRanges.Add(@($rng.Start, $rng.End))
}

# ---------------------------------------------------------------------------
# Step 3: every remaining run whose text is exactly the delimiter "{{" or
# "}}" gets bold + black color added (the template-tag marker runs, kept
# separate from the variable-name run between them).
# ---------------------------------------------------------------------------
foreach ($delim in @("{{", "}}")) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    while ($rng.Find.Execute($delim, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        $start = $rng.Start
        $end = $rng.End
        if (-not (InsideSkipRange $start)) {
            $rng.Font.Bold = $true
            $rng.Font.Color = 0
        }
        $rng.Collapse(0)
    }
}
